# Update NATMI LR-pair output (Il6-Il6st) with recomputed values based on new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2022703333333333
$ws.Range("H2").Value = 0.606811
$ws.Range("I2").Value = 0.01168815774551004
$ws.Range("J2").Value = 0.01168815774551004
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 4.010552858719
$ws.Range("R2").Value = 36.094975728471
$ws.Range("S2").Value = 0.001798093416304884
$ws.Range("T2").Value = 0.001798093416304884

# Row 3
$ws.Range("G3").Value = 0.2022703333333333
$ws.Range("H3").Value = 0.606811
$ws.Range("I3").Value = 0.01168815774551004
$ws.Range("J3").Value = 0.01168815774551004
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("Q3").Value = 17.21717209017367
$ws.Range("R3").Value = 154.954548811563
$ws.Range("S3").Value = 0.007719156154600022
$ws.Range("T3").Value = 0.007719156154600019

# Row 4
$ws.Range("G4").Value = 0.2022703333333333
$ws.Range("H4").Value = 0.606811
$ws.Range("I4").Value = 0.01168815774551004
$ws.Range("J4").Value = 0.01168815774551004
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 4.842096582262777
$ws.Range("R4").Value = 43.578869240365
$ws.Range("S4").Value = 0.002170908174605138
$ws.Range("T4").Value = 0.002170908174605137

# Row 5
$ws.Range("I5").Value = 0.8358439174604506
$ws.Range("J5").Value = 0.8358439174604506
$ws.Range("M5").Value = 19.827687
$ws.Range("N5").Value = 59.483061
$ws.Range("O5").Value = 0.1538389073329896
$ws.Range("P5").Value = 0.1538389073329896
$ws.Range("Q5").Value = 286.802786683952
$ws.Range("R5").Value = 2581.225080155568
$ws.Range("S5").Value = 0.1285853149630413
$ws.Range("T5").Value = 0.1285853149630413

# Row 6
$ws.Range("I6").Value = 0.8358439174604506
$ws.Range("J6").Value = 0.8358439174604506
$ws.Range("O6").Value = 0.6604253914664442
$ws.Range("P6").Value = 0.6604253914664441
$ws.Range("S6").Value = 0.5520125463936644
$ws.Range("T6").Value = 0.5520125463936643

# Row 7
$ws.Range("I7").Value = 0.8358439174604506
$ws.Range("J7").Value = 0.8358439174604506
$ws.Range("M7").Value = 23.93873833333333
$ws.Range("N7").Value = 71.816215
$ws.Range("O7").Value = 0.1857357012005663
$ws.Range("P7").Value = 0.1857357012005663
$ws.Range("Q7").Value = 346.2681685311022
$ws.Range("R7").Value = 3116.41351677992
$ws.Range("S7").Value = 0.1552460561037451
$ws.Range("T7").Value = 0.155246056103745

# Row 8
$ws.Range("G8").Value = 2.638545666666667
$ws.Range("H8").Value = 7.915637
$ws.Range("I8").Value = 0.1524679247940394
$ws.Range("J8").Value = 0.1524679247940394
$ws.Range("M8").Value = 19.827687
$ws.Range("N8").Value = 59.483061
$ws.Range("O8").Value = 0.1538389073329896
$ws.Range("P8").Value = 0.1538389073329896
$ws.Range("Q8").Value = 52.31625761387301
$ws.Range("R8").Value = 470.846318524857
$ws.Range("S8").Value = 0.02345549895364347
$ws.Range("T8").Value = 0.02345549895364346

# Row 9
$ws.Range("G9").Value = 2.638545666666667
$ws.Range("H9").Value = 7.915637
$ws.Range("I9").Value = 0.1524679247940394
$ws.Range("J9").Value = 0.1524679247940394
$ws.Range("O9").Value = 0.6604253914664442
$ws.Range("P9").Value = 0.6604253914664441
$ws.Range("Q9").Value = 224.5919807524024
$ws.Range("R9").Value = 2021.327826771621
$ws.Range("S9").Value = 0.1006936889181799
$ws.Range("T9").Value = 0.1006936889181799

# Row 10
$ws.Range("G10").Value = 2.638545666666667
$ws.Range("H10").Value = 7.915637
$ws.Range("I10").Value = 0.1524679247940394
$ws.Range("J10").Value = 0.1524679247940394
$ws.Range("M10").Value = 23.93873833333333
$ws.Range("N10").Value = 71.816215
$ws.Range("O10").Value = 0.1857357012005663
$ws.Range("P10").Value = 0.1857357012005663
$ws.Range("Q10").Value = 63.16345429488389
$ws.Range("R10").Value = 568.471088653955
$ws.Range("S10").Value = 0.02831873692221613
$ws.Range("T10").Value = 0.02831873692221612
